# Update cached market-price / profit values across all Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1549.2916
$ws.Range("J17").Value = 1549.2916
$ws.Range("L17").Value = 4647.8748
$ws.Range("N17").Value = -4983.8748

$ws.Range("H40").Value = 2866.5
$ws.Range("I40").Value = 2533.3333
$ws.Range("J40").Value = 3199.6667
$ws.Range("K40").Value = 2533.3333
$ws.Range("L40").Value = 3199.6667
$ws.Range("M40").Value = -2358.3333
$ws.Range("N40").Value = -3549.6667

$ws.Range("H93").Value = 34988
$ws.Range("J93").Value = 34988
$ws.Range("L93").Value = 34988
$ws.Range("N93").Value = -39980

$ws.Range("H98").Value = 1302.9333
$ws.Range("J98").Value = 3074.4
$ws.Range("L98").Value = 3074.4
$ws.Range("N98").Value = -6070.4

$ws.Range("H100").Value = 930.5
$ws.Range("I100").Value = 896.5
$ws.Range("J100").Value = 998.5
$ws.Range("K100").Value = 896.5
$ws.Range("L100").Value = 998.5
$ws.Range("M100").Value = -355.5
$ws.Range("N100").Value = -2080.5

$ws.Range("H107").Value = 1103.5333
$ws.Range("I107").Value = 252.16667
$ws.Range("K107").Value = 252.16667
$ws.Range("M107").Value = 1667.83333

$ws.Range("H122").Value = 1302.9333
$ws.Range("J122").Value = 3074.4
$ws.Range("L122").Value = 9223.200000000001
$ws.Range("N122").Value = -14123.2

$ws.Range("H126").Value = 100499.5
$ws.Range("J126").Value = 100499.5
$ws.Range("L126").Value = 100499.5
$ws.Range("N126").Value = -110379.5

$ws.Range("H132").Value = 335171
$ws.Range("I132").Value = 1535.6086
$ws.Range("K132").Value = 4606.825800000001
$ws.Range("M132").Value = -2076.825800000001

$ws.Range("H138").Value = 2747.9756
$ws.Range("I138").Value = 3547.9092
$ws.Range("J138").Value = 2454.6667
$ws.Range("K138").Value = 10643.7276
$ws.Range("L138").Value = 7364.000100000001
$ws.Range("M138").Value = -5503.7276
$ws.Range("N138").Value = -17644.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1338.8
$ws.Range("J26").Value = 580
$ws.Range("L26").Value = 580
$ws.Range("N26").Value = -1240

$ws.Range("H32").Value = 2589.2034
$ws.Range("I32").Value = 1280.0878
$ws.Range("K32").Value = 1280.0878
$ws.Range("M32").Value = -993.0878

$ws.Range("H45").Value = 1358.7142
$ws.Range("J45").Value = 1264.6666
$ws.Range("L45").Value = 1264.6666
$ws.Range("N45").Value = -2018.6666

$ws.Range("H61").Value = 2822.9614
$ws.Range("I61").Value = 2955.9443
$ws.Range("K61").Value = 2955.9443
$ws.Range("M61").Value = -2743.9443

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H97").Value = 1043.381
$ws.Range("I97").Value = 595.55
$ws.Range("K97").Value = 595.55
$ws.Range("M97").Value = -99.54999999999995

$ws.Range("H129").Value = 20990
$ws.Range("I129").Value = 20990
$ws.Range("K129").Value = 20990
$ws.Range("M129").Value = -15990

$ws.Range("H132").Value = 2303.625
$ws.Range("I132").Value = 2222
$ws.Range("K132").Value = 6666
$ws.Range("M132").Value = -4136

$ws.Range("H136").Value = 2822.9614
$ws.Range("I136").Value = 2955.9443
$ws.Range("K136").Value = 8867.832900000001
$ws.Range("M136").Value = -6317.832900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5299.4287
$ws.Range("I99").Value = 5281.4
$ws.Range("J99").Value = 5344.5
$ws.Range("K99").Value = 5281.4
$ws.Range("L99").Value = 5344.5
$ws.Range("M99").Value = -3783.4
$ws.Range("N99").Value = -8340.5

$ws.Range("H105").Value = 2238.8572
$ws.Range("I105").Value = 1449.8889
$ws.Range("K105").Value = 1449.8889
$ws.Range("M105").Value = 297.1111000000001

$ws.Range("H134").Value = 1668.5
$ws.Range("I134").Value = 1422.3
$ws.Range("J134").Value = 2899.5
$ws.Range("K134").Value = 4266.9
$ws.Range("L134").Value = 8698.5
$ws.Range("M134").Value = -1731.9
$ws.Range("N134").Value = -13768.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1685.0513
$ws.Range("I31").Value = 1238.5454
$ws.Range("J31").Value = 4140.8335
$ws.Range("K31").Value = 1238.5454
$ws.Range("L31").Value = 4140.8335
$ws.Range("M31").Value = -943.5454
$ws.Range("N31").Value = -4730.8335

$ws.Range("H34").Value = 1685.0513
$ws.Range("I34").Value = 1238.5454
$ws.Range("J34").Value = 4140.8335
$ws.Range("K34").Value = 1238.5454
$ws.Range("L34").Value = 4140.8335
$ws.Range("M34").Value = -1036.5454
$ws.Range("N34").Value = -4544.8335

$ws.Range("H41").Value = 11799.9
$ws.Range("J41").Value = 12499.875
$ws.Range("L41").Value = 12499.875
$ws.Range("N41").Value = -13355.875

$ws.Range("H58").Value = 3374.25
$ws.Range("I58").Value = 3374.25
$ws.Range("K58").Value = 3374.25
$ws.Range("M58").Value = -3171.25

$ws.Range("H93").Value = 18288.334
$ws.Range("I93").Value = 7085.143
$ws.Range("K93").Value = 7085.143
$ws.Range("M93").Value = -5213.143

$ws.Range("H99").Value = 2738.6667
$ws.Range("J99").Value = 2802
$ws.Range("L99").Value = 2802
$ws.Range("N99").Value = -5798

$ws.Range("H122").Value = 2362
$ws.Range("J122").Value = 3998
$ws.Range("L122").Value = 11994
$ws.Range("N122").Value = -16894

$ws.Range("H126").Value = 2738.6667
$ws.Range("J126").Value = 2802
$ws.Range("L126").Value = 8406
$ws.Range("N126").Value = -13346

$ws.Range("H132").Value = 1749.75
$ws.Range("J132").Value = 1999.5
$ws.Range("L132").Value = 5998.5
$ws.Range("N132").Value = -11058.5

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 3374.25
$ws.Range("I136").Value = 3374.25
$ws.Range("K136").Value = 10122.75
$ws.Range("M136").Value = -7572.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 700.5
$ws.Range("I7").Value = 700.5
$ws.Range("K7").Value = 2101.5
$ws.Range("M7").Value = -1989.5

$ws.Range("H34").Value = 244.8
$ws.Range("I34").Value = 244.8
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 734.4000000000001
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -650.4000000000001
$ws.Range("N34").ClearContents()

$ws.Range("H39").Value = 2453.8572
$ws.Range("J39").Value = 2453.8572
$ws.Range("L39").Value = 7361.571599999999
$ws.Range("N39").Value = -7949.571599999999

$ws.Range("H40").Value = 39.909092
$ws.Range("I40").Value = 37.666668
$ws.Range("J40").Value = 50
$ws.Range("K40").Value = 150.666672
$ws.Range("L40").Value = 200
$ws.Range("M40").Value = -81.66667200000001
$ws.Range("N40").Value = -338

$ws.Range("H55").Value = 3999.6667
$ws.Range("J55").Value = 4000
$ws.Range("L55").Value = 12000
$ws.Range("N55").Value = -12354

$ws.Range("H58").Value = 20000
$ws.Range("I58").Value = 20000
$ws.Range("K58").Value = 60000
$ws.Range("M58").Value = -59872

$ws.Range("H113").Value = 845.8
$ws.Range("I113").Value = 628.5
$ws.Range("J113").Value = 990.6667
$ws.Range("K113").Value = 1885.5
$ws.Range("L113").Value = 2972.0001
$ws.Range("M113").Value = 284.5
$ws.Range("N113").Value = -7312.0001

$ws.Range("H129").Value = 3235.353
$ws.Range("I129").Value = 607.8889
$ws.Range("K129").Value = 1823.6667
$ws.Range("M129").Value = 3176.3333

$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 45000
$ws.Range("M132").Value = -42470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1500.4736
$ws.Range("I97").Value = 1238.4375
$ws.Range("J97").Value = 2898
$ws.Range("K97").Value = 1238.4375
$ws.Range("L97").Value = 2898
$ws.Range("M97").Value = -742.4375
$ws.Range("N97").Value = -3890

$ws.Range("H132").Value = 2753.0667
$ws.Range("I132").Value = 2709.182
$ws.Range("K132").Value = 8127.545999999999
$ws.Range("M132").Value = -5597.545999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 1588.6666
$ws.Range("I46").Value = 1070.5
$ws.Range("J46").Value = 2625
$ws.Range("K46").Value = 1070.5
$ws.Range("L46").Value = 2625
$ws.Range("M46").Value = -882.5
$ws.Range("N46").Value = -3001

$ws.Range("H127").Value = 67500
$ws.Range("J127").Value = 67500
$ws.Range("L127").Value = 67500
$ws.Range("N127").Value = -77420

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 7500
$ws.Range("J7").Value = 7500
$ws.Range("L7").Value = 7500
$ws.Range("N7").Value = -7726

$ws.Range("H132").Value = 1658.25
$ws.Range("I132").Value = 1590.2
$ws.Range("K132").Value = 4770.6
$ws.Range("M132").Value = -2240.6

$ws.Range("H136").Value = 3398.6667
$ws.Range("J136").Value = 443
$ws.Range("L136").Value = 1329
$ws.Range("N136").Value = -6429
